# Update BOC USD rates (auto)
# Adds the 2026-01-03 00:00:05 publish (BOC USD) to:
#   - "All Published Values": new row 20
#   - "Daily Summary": new row in the "Day Averages" block (row 5) and a
#     new row in the "Day First Published" block (row 11, after the
#     insert shifts the existing rows 5-9 down to 6-10)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value into a cell as literal TEXT (no Excel smart
# number/date inference), without leaving a lingering custom style on
# the cell - mirrors the plain inlineStr cells the source file uses.
# ---------------------------------------------------------------------
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

# =======================================================================
# Sheet 1: "All Published Values" - append row 20
# =======================================================================
$ws1 = $wb.Worksheets.Item("All Published Values")

Set-TextValue $ws1.Range("A20") "2026-01-03"
Set-TextValue $ws1.Range("B20") "2026-01-03 00:00:05"
Set-TextValue $ws1.Range("C20") "697.85"
Set-TextValue $ws1.Range("D20") "697.85"
Set-TextValue $ws1.Range("E20") "700.79"
Set-TextValue $ws1.Range("F20") "700.79"
Set-TextValue $ws1.Range("G20") "702.88"
Set-TextValue $ws1.Range("H20") "2026/01/03 00:00:05"
Set-TextValue $ws1.Range("I20") "2026-01-02 16:26:26"
Set-TextValue $ws1.Range("J20") "https://www.bankofchina.com/sourcedb/whpj/enindex_1619.html"

# Re-extend the AutoFilter over the new last row.
$ws1.AutoFilterMode = $false
$null = $ws1.Range("A1:J20").AutoFilter()

# Keep the hidden _FilterDatabase defined name lined up with the filter.
$names1 = $wb.Names
for ($i = 1; $i -le $names1.Count(); $i++) {
    $n = $names1.Item($i)
    if ($n.Name() -eq "All Published Values!_FilterDatabase") {
        $n.RefersTo = "='All Published Values'!" + '$A$1:$J$20'
    }
}

# =======================================================================
# Sheet 2: "Daily Summary"
# =======================================================================
$ws2 = $wb.Worksheets.Item("Daily Summary")

# Insert a new row 5 - pushes the old row5..row9 ("" / header / blank /
# header / 2026-01-02 first-published row) down to row6..row10.
$ws2.Rows.Item(5).Insert()

# New "Day Averages (Middle)" row for 2026-01-03.
Set-TextValue $ws2.Range("A5") "2026-01-03"
$ws2.Range("B5").Value = 1
$ws2.Range("C5").Value = 702.88
$ws2.Range("D5").Value = 702.88
$ws2.Range("E5").Value = 702.88

# New "Day First Published (Middle)" row for 2026-01-03 (row 11, right
# after the shifted 2026-01-02 row now sitting at row 10).
Set-TextValue $ws2.Range("A11") "2026-01-03"
Set-TextValue $ws2.Range("B11") "2026-01-03 00:00:05"
Set-TextValue $ws2.Range("C11") "702.88"
Set-TextValue $ws2.Range("D11") "2026/01/03 00:00:05"

# Re-extend the AutoFilter over the new last row of the top block.
$ws2.AutoFilterMode = $false
$null = $ws2.Range("A1:E5").AutoFilter()

$names2 = $wb.Names
for ($i = 1; $i -le $names2.Count(); $i++) {
    $n = $names2.Item($i)
    if ($n.Name() -eq "Daily Summary!_FilterDatabase") {
        $n.RefersTo = "='Daily Summary'!" + '$A$1:$E$5'
    }
}

Write-Host "done"
